$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting rows 51-95 down to 52-96.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new weekly price record.
$ws.Cells.Item(51, 1).Value = 8
$ws.Cells.Item(51, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(51, 3).Value = "Coquimbo"
$ws.Cells.Item(51, 4).Value = 44729
$ws.Cells.Item(51, 5).Value = 4
$ws.Cells.Item(51, 6).Value = 100112052
$ws.Cells.Item(51, 7).Value = "Albahaca"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 1140
$ws.Cells.Item(51, 11).Value = 3500
$ws.Cells.Item(51, 12).Value = 4000
$ws.Cells.Item(51, 13).Value = 3750
$ws.Cells.Item(51, 14).Value = "$/paquete"
$ws.Cells.Item(51, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(51, 16).Value = 3750
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = "Hortaliza"
